$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '64.591.49'
$ws.Range('E2').Value = '  +0.81%  '

$ws.Range('D3').Value = '3.151.68'
$ws.Range('E3').Value = '  +0.25%  '

$ws.Range('E4').Value = '  +0.02%  '

Set-TextValue $ws 'D5' '608.47'
$ws.Range('E5').Value = '  -0.66%  '

Set-TextValue $ws 'D6' '144.11'
$ws.Range('E6').Value = '  -1.52%  '

$ws.Range('E7').Value = '  -0.12%  '

$ws.Range('D8').Value = '3.148.38'
$ws.Range('E8').Value = '  +0.51%  '

Set-TextValue $ws 'D9' '0.524'
$ws.Range('E9').Value = '  +0.14%  '

$ws.Range('E10').Value = '  +0.64%  '

Set-TextValue $ws 'D11' '5.43'
$ws.Range('E11').Value = '  +2.40%  '

Set-TextValue $ws 'D12' '0.469'
$ws.Range('E12').Value = '  -0.40%  '

$ws.Range('E13').Value = '  +3.71%  '

Set-TextValue $ws 'D14' '35.45'
$ws.Range('E14').Value = '  +0.48%  '

$ws.Range('D15').Value = '3.669.08'
$ws.Range('E15').Value = '  +0.38%  '

$ws.Range('E16').Value = '  +2.60%  '

$ws.Range('D17').Value = '64.429.47'
$ws.Range('E17').Value = '  +0.53%  '

$ws.Range('D18').Value = '3.151.86'
$ws.Range('E18').Value = '  +0.45%  '

$ws.Range('E19').Value = '  +0.50%  '

Set-TextValue $ws 'D20' '481.37'
$ws.Range('E20').Value = '  +1.07%  '

Set-TextValue $ws 'D21' '14.66'
$ws.Range('E21').Value = '  -0.17%  '

Set-TextValue $ws 'D22' '0.717'
$ws.Range('E22').Value = '  +1.99%  '

Set-TextValue $ws 'D23' '7.71'
$ws.Range('E23').Value = '  -0.44%  '

Set-TextValue $ws 'D24' '85.31'
$ws.Range('E24').Value = '  +2.14%  '

Set-TextValue $ws 'D25' '13.45'
$ws.Range('E25').Value = '  -0.90%  '

$ws.Range('E26').Value = '  +0.11%  '

$ws.Range('E27').Value = '  -1.30%  '

Set-TextValue $ws 'D28' '8.44'
$ws.Range('E28').Value = '  +0.45%  '

Set-TextValue $ws 'D29' '7.20'
$ws.Range('E29').Value = '  +7.08%  '

$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D30' '2.05'
$ws.Range('E30').Value = '  -5.27%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D31' '0.115'
$ws.Range('E31').Value = '  +1.82%  '

Set-TextValue $ws 'D32' '27.05'
$ws.Range('E32').Value = '  +3.44%  '

$ws.Range('E33').Value = '  +0.01%  '

Set-TextValue $ws 'D34' '2.66'
$ws.Range('E34').Value = '  -2.49%  '

$ws.Range('E35').Value = '  -1.47%  '

Set-TextValue $ws 'D36' '5.99'
$ws.Range('E36').Value = '  +0.75%  '

$ws.Range('D37').Value = '0.0₃0770'
$ws.Range('E37').Value = '  +5.34%  '

Set-TextValue $ws 'D38' '52.46'
$ws.Range('E38').Value = '  -1.76%  '

Set-TextValue $ws 'D39' '3.04'
$ws.Range('E39').Value = '  +4.54%  '

Set-TextValue $ws 'D40' '447.48'
$ws.Range('E40').Value = '  -3.22%  '

$ws.Range('E41').Value = '  +0.45%  '

$ws.Range('E42').Value = '  +1.12%  '

Set-TextValue $ws 'D43' '8.25'
$ws.Range('E43').Value = '  -1.56%  '

$ws.Range('D44').Value = '2.871.77'
$ws.Range('E44').Value = '  +1.00%  '

Set-TextValue $ws 'D45' '0.263'
$ws.Range('E45').Value = '  -0.96%  '

Set-TextValue $ws 'D46' '2.24'
$ws.Range('E46').Value = '  -0.39%  '

Set-TextValue $ws 'D47' '2.41'
$ws.Range('E47').Value = '  +2.12%  '

$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws 'D48' '0.999'
$ws.Range('E48').Value = '  -0.01%  '

$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D49' '26.24'
$ws.Range('E49').Value = '  -0.15%  '

$ws.Range('E50').Value = '  -0.20%  '

Set-TextValue $ws 'D51' '119.61'
$ws.Range('E51').Value = '  +1.08%  '
